# Auto-generated edit script applying value updates to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 604.7778
$ws.Range("I32").Value = 450.5
$ws.Range("J32").Value = 648.8570999999999
$ws.Range("K32").Value = 450.5
$ws.Range("L32").Value = 648.8570999999999
$ws.Range("M32").Value = -124.5
$ws.Range("N32").Value = -1300.8571
$ws.Range("H55").Value = 160.38461
$ws.Range("I55").Value = 134.28572
$ws.Range("J55").Value = 190.83333
$ws.Range("K55").Value = 134.28572
$ws.Range("L55").Value = 190.83333
$ws.Range("M55").Value = 79.71428
$ws.Range("N55").Value = -618.8333299999999
$ws.Range("H112").Value = 7576983
$ws.Range("I112").Value = 699.6667
$ws.Range("J112").Value = 8265736
$ws.Range("K112").Value = 2099.0001
$ws.Range("L112").Value = 24797208
$ws.Range("M112").Value = -991.0001000000002
$ws.Range("N112").Value = -24799424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N2").ClearContents()
$ws.Range("H2").Value = 36649.18
$ws.Range("I2").Value = 36649.18
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 36649.18
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = -36536.18
$ws.Range("H32").Value = 20254.12
$ws.Range("I32").Value = 4444.86
$ws.Range("J32").Value = 67681.89999999999
$ws.Range("K32").Value = 4444.86
$ws.Range("L32").Value = 67681.89999999999
$ws.Range("M32").Value = -4157.86
$ws.Range("N32").Value = -68255.89999999999
$ws.Range("H45").Value = 889.2857
$ws.Range("I45").Value = 861.1111
$ws.Range("J45").Value = 940
$ws.Range("K45").Value = 861.1111
$ws.Range("L45").Value = 940
$ws.Range("M45").Value = -484.1111
$ws.Range("N45").Value = -1694
$ws.Range("H61").Value = 3330.923
$ws.Range("I61").Value = 2246
$ws.Range("J61").Value = 6947.3335
$ws.Range("K61").Value = 2246
$ws.Range("L61").Value = 6947.3335
$ws.Range("M61").Value = -2034
$ws.Range("N61").Value = -7371.3335
$ws.Range("H63").Value = 9465.950000000001
$ws.Range("I63").Value = 11327.571
$ws.Range("J63").Value = 5122.1665
$ws.Range("K63").Value = 11327.571
$ws.Range("L63").Value = 5122.1665
$ws.Range("M63").Value = -10641.571
$ws.Range("N63").Value = -6494.1665
$ws.Range("H66").Value = 9465.950000000001
$ws.Range("I66").Value = 11327.571
$ws.Range("J66").Value = 5122.1665
$ws.Range("K66").Value = 56637.855
$ws.Range("L66").Value = 25610.8325
$ws.Range("M66").Value = -53205.855
$ws.Range("N66").Value = -32474.8325
$ws.Range("N116").ClearContents()
$ws.Range("H116").Value = 36649.18
$ws.Range("I116").Value = 36649.18
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 36649.18
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = -34355.18
$ws.Range("H136").Value = 3330.923
$ws.Range("I136").Value = 2246
$ws.Range("J136").Value = 6947.3335
$ws.Range("K136").Value = 6738
$ws.Range("L136").Value = 20842.0005
$ws.Range("M136").Value = -4188
$ws.Range("N136").Value = -25942.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N3").ClearContents()
$ws.Range("H3").Value = 36649.18
$ws.Range("I3").Value = 36649.18
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 36649.18
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = -36535.18
$ws.Range("H20").Value = 1401.8948
$ws.Range("I20").Value = 1103.0834
$ws.Range("J20").Value = 1914.1428
$ws.Range("K20").Value = 1103.0834
$ws.Range("L20").Value = 1914.1428
$ws.Range("M20").Value = -856.0834
$ws.Range("N20").Value = -2408.1428
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("N86").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("N89").Value = 0
$ws.Range("H102").Value = 8217.714
$ws.Range("I102").Value = 8217.714
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 8217.714
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -4972.714
$ws.Range("H107").Value = 681.93335
$ws.Range("I107").Value = 725.3077
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 725.3077
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1194.6923
$ws.Range("N107").Value = -4240
$ws.Range("H134").Value = 4191.4414
$ws.Range("I134").Value = 2398.3635
$ws.Range("J134").Value = 7478.75
$ws.Range("K134").Value = 7195.0905
$ws.Range("L134").Value = 22436.25
$ws.Range("M134").Value = -4660.0905
$ws.Range("N134").Value = -27506.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 69.75
$ws.Range("I7").Value = 76.333336
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 76.333336
$ws.Range("L7").Value = 50
$ws.Range("M7").Value = 36.666664
$ws.Range("N7").Value = -276
$ws.Range("H58").Value = 2619.375
$ws.Range("I58").Value = 1552.7273
$ws.Range("J58").Value = 4966
$ws.Range("K58").Value = 1552.7273
$ws.Range("L58").Value = 4966
$ws.Range("M58").Value = -1349.7273
$ws.Range("N58").Value = -5372
$ws.Range("H107").Value = 362.30768
$ws.Range("I107").Value = 291.72726
$ws.Range("J107").Value = 750.5
$ws.Range("K107").Value = 291.72726
$ws.Range("L107").Value = 750.5
$ws.Range("M107").Value = 1628.27274
$ws.Range("N107").Value = -4590.5
$ws.Range("H132").Value = 1962.5555
$ws.Range("I132").Value = 1151
$ws.Range("J132").Value = 4072.6
$ws.Range("K132").Value = 3453
$ws.Range("L132").Value = 12217.8
$ws.Range("M132").Value = -923
$ws.Range("N132").Value = -17277.8
$ws.Range("H136").Value = 2619.375
$ws.Range("I136").Value = 1552.7273
$ws.Range("J136").Value = 4966
$ws.Range("K136").Value = 4658.1819
$ws.Range("L136").Value = 14898
$ws.Range("M136").Value = -2108.1819
$ws.Range("N136").Value = -19998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 821
$ws.Range("I23").Value = 1194.7778
$ws.Range("J23").Value = 340.42856
$ws.Range("K23").Value = 3584.3334
$ws.Range("L23").Value = 1021.28568
$ws.Range("M23").Value = -3349.3334
$ws.Range("N23").Value = -1491.28568
$ws.Range("H113").Value = 512.4167
$ws.Range("I113").Value = 475
$ws.Range("J113").Value = 517.09375
$ws.Range("K113").Value = 1425
$ws.Range("L113").Value = 1551.28125
$ws.Range("M113").Value = 745
$ws.Range("N113").Value = -5891.28125
$ws.Range("H122").Value = 733.5263
$ws.Range("I122").Value = 327.66666
$ws.Range("J122").Value = 1098.8
$ws.Range("K122").Value = 2948.99994
$ws.Range("L122").Value = 9889.199999999999
$ws.Range("M122").Value = -498.9999399999997
$ws.Range("N122").Value = -14789.2
$ws.Range("H131").Value = 9010659
$ws.Range("I131").Value = 540
$ws.Range("J131").Value = 9525523
$ws.Range("K131").Value = 1620
$ws.Range("L131").Value = 28576569
$ws.Range("M131").Value = 3420
$ws.Range("N131").Value = -28586649
$ws.Range("H137").Value = 6315895
$ws.Range("I137").Value = 10001372
$ws.Range("J137").Value = 173433
$ws.Range("K137").Value = 30004116
$ws.Range("L137").Value = 520299
$ws.Range("M137").Value = -29999016
$ws.Range("N137").Value = -530499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5049.8203
$ws.Range("I70").Value = 5374.276
$ws.Range("J70").Value = 4108.9
$ws.Range("K70").Value = 5374.276
$ws.Range("L70").Value = 4108.9
$ws.Range("M70").Value = -5104.276
$ws.Range("N70").Value = -4648.9
$ws.Range("H73").Value = 5049.8203
$ws.Range("I73").Value = 5374.276
$ws.Range("J73").Value = 4108.9
$ws.Range("K73").Value = 5374.276
$ws.Range("L73").Value = 4108.9
$ws.Range("M73").Value = -4438.276
$ws.Range("N73").Value = -5980.9
$ws.Range("M74").ClearContents()
$ws.Range("H74").Value = 49933
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 49933
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = 49933
$ws.Range("N74").Value = -51805
$ws.Range("M77").ClearContents()
$ws.Range("H77").Value = 49933
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 49933
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = 149799
$ws.Range("N77").Value = -159159
$ws.Range("H132").Value = 2716.9736
$ws.Range("I132").Value = 2159.303
$ws.Range("J132").Value = 6397.6
$ws.Range("K132").Value = 6477.909
$ws.Range("L132").Value = 19192.8
$ws.Range("M132").Value = -3947.909
$ws.Range("N132").Value = -24252.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11415
$ws.Range("I22").Value = 1250
$ws.Range("J22").Value = 12544.444
$ws.Range("K22").Value = 1250
$ws.Range("L22").Value = 12544.444
$ws.Range("M22").Value = -955
$ws.Range("N22").Value = -13134.444
$ws.Range("H27").Value = 11415
$ws.Range("I27").Value = 1250
$ws.Range("J27").Value = 12544.444
$ws.Range("K27").Value = 1250
$ws.Range("L27").Value = 12544.444
$ws.Range("M27").Value = -1143
$ws.Range("N27").Value = -12758.444
$ws.Range("H46").Value = 1375
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -1876
$ws.Range("H55").Value = 501.33334
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 465.7143
$ws.Range("K55").Value = 1000
$ws.Range("L55").Value = 465.7143
$ws.Range("M55").Value = -827
$ws.Range("N55").Value = -811.7143
$ws.Range("H68").Value = 1948.625
$ws.Range("I68").Value = 1912.7142
$ws.Range("J68").Value = 2200
$ws.Range("K68").Value = 1912.7142
$ws.Range("L68").Value = 2200
$ws.Range("M68").Value = -1163.7142
$ws.Range("N68").Value = -3698
$ws.Range("H71").Value = 1948.625
$ws.Range("I71").Value = 1912.7142
$ws.Range("J71").Value = 2200
$ws.Range("K71").Value = 9563.571
$ws.Range("L71").Value = 11000
$ws.Range("M71").Value = -5819.571
$ws.Range("N71").Value = -18488
$ws.Range("H100").Value = 1955898.8
$ws.Range("I100").Value = 4809866
$ws.Range("J100").Value = 3184.2104
$ws.Range("K100").Value = 4809866
$ws.Range("L100").Value = 3184.2104
$ws.Range("M100").Value = -4809325
$ws.Range("N100").Value = -4266.2104
$ws.Range("H132").Value = 3831.1892
$ws.Range("I132").Value = 2358.8333
$ws.Range("J132").Value = 6549.385
$ws.Range("K132").Value = 7076.499899999999
$ws.Range("L132").Value = 19648.155
$ws.Range("M132").Value = -4546.499899999999
$ws.Range("N132").Value = -24708.155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3068.5112
$ws.Range("I132").Value = 3362.9375
$ws.Range("J132").Value = 2343.7693
$ws.Range("K132").Value = 10088.8125
$ws.Range("L132").Value = 7031.3079
$ws.Range("M132").Value = -7558.8125
$ws.Range("N132").Value = -12091.3079
$ws.Range("H136").Value = 1644.3334
$ws.Range("I136").Value = 914.5454999999999
$ws.Range("J136").Value = 4855.4
$ws.Range("K136").Value = 2743.6365
$ws.Range("L136").Value = 14566.2
$ws.Range("M136").Value = -193.6364999999996
$ws.Range("N136").Value = -19666.2
